# Rearranges rows 8-23 on the active sheet: each destination row receives the
# full content (all columns A:AY) that a particular source row held before the
# edit. This reproduces an upstream re-sort/shuffle of the same 16 records
# without altering their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 8
$lastRow = 23
$firstCol = 1
$lastCol = 51

# destination row -> source row (source row's pre-edit content moves to the
# destination row)
$mapping = @{
    8  = 20
    9  = 10
    10 = 16
    11 = 12
    12 = 9
    13 = 19
    14 = 21
    15 = 22
    16 = 15
    17 = 8
    18 = 11
    19 = 17
    20 = 13
    21 = 23
    22 = 18
    23 = 14
}

# 1) Snapshot every cell's value (and presence - $null means the cell didn't
#    exist at all, as opposed to an empty string) before touching anything.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowVals
}

# 2) Clear every cell in the range so rows that should end up without a given
#    column (e.g. a column that didn't exist in the source row) truly have no
#    cell left behind.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r, $c).ClearContents()
    }
}

# 3) Write each destination row's cells back from its mapped source row's
#    snapshot, skipping cells that were absent (null) in the source. Values
#    that look like plain "YYYY-MM-DD" text (e.g. the Startdatum/Slutdatum
#    columns, stored as literal text in the source file) must be written as
#    text explicitly - otherwise Excel's normal typed-value auto-detection
#    would silently turn them into real date serials, which is not what the
#    source workbook had.
foreach ($dst in $mapping.Keys) {
    $src = $mapping[$dst]
    $rowVals = $snapshot[$src]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $v = $rowVals[$c]
        if ($v -ne $null) {
            $cell = $ws.Cells.Item($dst, $c)
            if (($v -is [string]) -and ($v -match '^\d{4}-\d{2}-\d{2}$')) {
                $cell.NumberFormat = "@"
            }
            $cell.Value = $v
        }
    }
}
